$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 already exists (B27 "Day 24", D27 date) -- just needs C27 filled in
# with the same wrapped narrative style used by C24:C26.
$ws.Range("C24").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null

# Build out rows 28-35 (Days 25-32). Column B/D formats come from the row
# directly above (plain text style + date style); column C format depends on
# whether the note is a single short line (plain style, like B) or a longer
# wrapped note (same wrap style as C24).
for ($r = 28; $r -le 34; $r++) {
    $prev = $r - 1
    $ws.Range("B$prev").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("D$prev").Copy() | Out-Null
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
}

$ws.Range("B34").Copy() | Out-Null
$ws.Range("B35").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# Column C styles: wrapped narrative style for the multi-line notes, plain
# style (matching column B) for the short single-line ones.
$ws.Range("C24").Copy() | Out-Null
$ws.Range("C27:C29").PasteSpecial(-4122) | Out-Null
$ws.Range("C31:C33").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("B30").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$ws.Range("C34:C35").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Values (filled in row-major order, matching how the shared-string
# table was originally populated: note text, then "Day N", for each row) ---
$ws.Range("C27").Value = " - Finished project 02 (Show me the Data Structure)`n- Worked on explanations for the project 02 problems"

$ws.Range("B28").Value = "Day 25"
$ws.Range("C28").Value = " - Fix problems on project 02 (Show me the Data Structure)`n- Search for books on the subject"
$ws.Range("D28").Value = 44116

$ws.Range("B29").Value = "Day 26"
$ws.Range("C29").Value = " - Learning about Binary Search`n- Reading the book ""Introduction to Algorithm"""
$ws.Range("D29").Value = 44117

$ws.Range("B30").Value = "Day 27"
$ws.Range("C30").Value = " - Worked on Tries and Heaps."
$ws.Range("D30").Value = 44118

$ws.Range("B31").Value = "Day 28"
$ws.Range("C31").Value = " - Learned about self-balancing trees`n- Studied about Red-Black Trees"
$ws.Range("D31").Value = 44119

$ws.Range("B32").Value = "Day 29"
$ws.Range("C32").Value = " - Start studying sorting algorithms.`n- I've learned about Bubble Sorted.`n- I've solved some exercises."
$ws.Range("D32").Value = 44120

$ws.Range("C33").Value = " - I've learned about Merge Sort.`n- I've solved some exercises."
$ws.Range("B33").Value = "Day 30"
$ws.Range("D33").Value = 44121

$ws.Range("B34").Value = "Day 31"
$ws.Range("D34").Value = 44122

$ws.Range("B35").Value = "Day 32"

$ws.Range("C34").Value = "I've learned about Quick Sort"
$ws.Range("C35").Value = "I've learned about Heap Sort"

$ws.Range("G25").Select() | Out-Null
